# Add two new impact items (combustor unit with energy recovery):
#   - WoodPellet
#   - StoneWool
# as new rows 21/22 on every sheet of the workbook (info, GWP, H_Ecosystems,
# H_Health, H_Resources).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) info sheet: ID / unit columns only
# ---------------------------------------------------------------------------
$info = $wb.Worksheets.Item("info")

$info.Range("A21").Value = "WoodPellet"
$info.Range("B21").Value = "kg"

$info.Range("A22").Value = "StoneWool"
$info.Range("B22").Value = "kg"

$info.Range("E22").Select()

# ---------------------------------------------------------------------------
# 2) GWP sheet: ID, unit, low/expected/high, distribution, references
# ---------------------------------------------------------------------------
$gwp = $wb.Worksheets.Item("GWP")

$gwp.Range("A21").Value = "WoodPellet"
$gwp.Range("B21").Value = "kg CO2-eq"
$gwp.Range("C21").Value = 0.14178389
$gwp.Range("D21").Formula = "=C21*0.9"
$gwp.Range("E21").Formula = "=C21*1.1"
$gwp.Range("F21").Value = "uniform"
$gwp.Range("G21").Value = "ecoinvent 3.8 - cutoff, TRACI"

$gwp.Range("A22").Value = "StoneWool"
$gwp.Range("B22").Value = "kg CO2-eq"
$gwp.Range("C22").Value = 1.2681741
$gwp.Range("D22").Formula = "=C22*0.9"
$gwp.Range("E22").Formula = "=C22*1.1"
$gwp.Range("F22").Value = "uniform"
$gwp.Range("G22").Value = "ecoinvent 3.8 - cutoff, TRACI"

$gwp.Activate()
$gwp.Range("F22:G22").Select()

# ---------------------------------------------------------------------------
# 3) H_Ecosystems sheet
# ---------------------------------------------------------------------------
$hEco = $wb.Worksheets.Item("H_Ecosystems")

$hEco.Range("A21").Value = "WoodPellet"
$hEco.Range("B21").Value = "points"
$hEco.Range("C21").Value = 0.043958951
$hEco.Range("D21").Formula = "=C21*0.9"
$hEco.Range("E21").Formula = "=C21*1.1"
$hEco.Range("F21").Value = "uniform"
$hEco.Range("G21").Value = "ecoinvent 3.8 - cutoff, ReCiPe Endpoint (H,A)"

$hEco.Range("A22").Value = "StoneWool"
$hEco.Range("B22").Value = "points"
$hEco.Range("C22").Value = 0.027597747
$hEco.Range("D22").Formula = "=C22*0.9"
$hEco.Range("E22").Formula = "=C22*1.1"
$hEco.Range("F22").Value = "uniform"
$hEco.Range("G22").Value = "ecoinvent 3.8 - cutoff, ReCiPe Endpoint (H,A)"

$hEco.Activate()
$excel.ActiveWindow.Zoom = 220
$hEco.Range("D24").Select()

# ---------------------------------------------------------------------------
# 4) H_Health sheet
# ---------------------------------------------------------------------------
$hHealth = $wb.Worksheets.Item("H_Health")

$hHealth.Range("A21").Value = "WoodPellet"
$hHealth.Range("B21").Value = "points"
$hHealth.Range("C21").Value = 0.0068199964
$hHealth.Range("D21").Formula = "=C21*0.9"
$hHealth.Range("E21").Formula = "=C21*1.1"
$hHealth.Range("F21").Value = "uniform"
$hHealth.Range("G21").Value = "ecoinvent 3.8 - cutoff, ReCiPe Endpoint (H,A)"

$hHealth.Range("A22").Value = "StoneWool"
$hHealth.Range("B22").Value = "points"
$hHealth.Range("C22").Value = 0.059336112
$hHealth.Range("D22").Formula = "=C22*0.9"
$hHealth.Range("E22").Formula = "=C22*1.1"
$hHealth.Range("F22").Value = "uniform"
$hHealth.Range("G22").Value = "ecoinvent 3.8 - cutoff, ReCiPe Endpoint (H,A)"

$hHealth.Activate()
$hHealth.Range("D23").Select()

# ---------------------------------------------------------------------------
# 5) H_Resources sheet
# ---------------------------------------------------------------------------
$hRes = $wb.Worksheets.Item("H_Resources")

$hRes.Range("A21").Value = "WoodPellet"
$hRes.Range("B21").Value = "points"
$hRes.Range("C21").Value = 0.006172004
$hRes.Range("D21").Formula = "=C21*0.9"
$hRes.Range("E21").Formula = "=C21*1.1"
$hRes.Range("F21").Value = "uniform"
$hRes.Range("G21").Value = "ecoinvent 3.8 - cutoff, ReCiPe Endpoint (H,A)"

$hRes.Range("A22").Value = "StoneWool"
$hRes.Range("B22").Value = "points"
$hRes.Range("C22").Value = 0.051289877
$hRes.Range("D22").Formula = "=C22*0.9"
$hRes.Range("E22").Formula = "=C22*1.1"
$hRes.Range("F22").Value = "uniform"
$hRes.Range("G22").Value = "ecoinvent 3.8 - cutoff, ReCiPe Endpoint (H,A)"

$hRes.Activate()
$hRes.Range("E26").Select()

# Leave the "info" sheet active/selected like the original workbook.
$info.Activate()
